$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume Number + Report Covering the Week dates) ---
$ws.Range("A8").Value = "Volume 30   Number  46"
$ws.Range("C9").Value = "Report Covering the Week  11/13/2023  Through  11/19/2023"

# --- Data table updates (rows 14-29) ---
# Set numeric/text values first, then fix up cell styles for cells whose
# underlying type changed (text "N/A" placeholder <-> real number), using
# stable donor cells (I14=integer style, K14=percent style, C14=text style)
# that are not touched by this edit, via Copy/PasteSpecial(xlPasteFormats).

# Row 14
$ws.Range("M14").Value = -76.470588235294

# Row 15
$ws.Range("D15").Value = 1
$ws.Range("I14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("G15").Value = 2
$ws.Range("J15").Value = 20
$ws.Range("K15").Value = -15
$ws.Range("N15").Value = -76.712328767123

# Row 16
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = -40
$ws.Range("I16").Value = 133
$ws.Range("J16").Value = 141
$ws.Range("K16").Value = -5.673758865248
$ws.Range("L16").Value = -15.822784810126
$ws.Range("M16").Value = -56.39344262295
$ws.Range("N16").Value = -87.941976427923

# Row 17
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 9
$ws.Range("F17").Value = 20
$ws.Range("G17").Value = 24
$ws.Range("H17").Value = -16.666666666666
$ws.Range("I17").Value = 265
$ws.Range("J17").Value = 299
$ws.Range("K17").Value = -11.371237458194
$ws.Range("L17").Value = -12.251655629139
$ws.Range("M17").Value = -8.93470790378
$ws.Range("N17").Value = -66.112531969309

# Row 18
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -70
$ws.Range("I18").Value = 125
$ws.Range("J18").Value = 188
$ws.Range("K18").Value = -33.510638297872
$ws.Range("L18").Value = -18.831168831168
$ws.Range("M18").Value = -41.860465116279
$ws.Range("N18").Value = -82.662968099861

# Row 19
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = -50
$ws.Range("F19").Value = 23
$ws.Range("G19").Value = 34
$ws.Range("H19").Value = -32.35294117647
$ws.Range("I19").Value = 301
$ws.Range("J19").Value = 387
$ws.Range("K19").Value = -22.222222222222
$ws.Range("L19").Value = -9.063444108761
$ws.Range("M19").Value = -4.746835443037
$ws.Range("N19").Value = -9.063444108761

# Row 20
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 200
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = -10
$ws.Range("I20").Value = 98
$ws.Range("J20").Value = 124
$ws.Range("K20").Value = -20.967741935483
$ws.Range("L20").Value = 7.692307692307
$ws.Range("M20").Value = 28.947368421052
$ws.Range("N20").Value = -81.885397412199

# Row 21
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = -13.636363636363
$ws.Range("G21").Value = 100
$ws.Range("H21").Value = -36
$ws.Range("I21").Value = 943
$ws.Range("J21").Value = 1164
$ws.Range("K21").Value = -18.986254295532
$ws.Range("L21").Value = -10.616113744075
$ws.Range("M21").Value = -24.074074074074
$ws.Range("N21").Value = -73.629753914988

# Row 22
$ws.Range("D22").Value = 1
$ws.Range("I14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("J22").Value = 14
$ws.Range("K22").Value = -21.428571428571

# Row 23
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 0
$ws.Range("I23").Value = 76
$ws.Range("J23").Value = 76
$ws.Range("L23").Value = -21.649484536082
$ws.Range("M23").Value = 2.702702702702

# Row 24
$ws.Range("C24").Value = 13
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = -7.142857142857
$ws.Range("F24").Value = 59
$ws.Range("G24").Value = 58
$ws.Range("H24").Value = 1.724137931034
$ws.Range("I24").Value = 757
$ws.Range("J24").Value = 786
$ws.Range("K24").Value = -3.689567430025
$ws.Range("L24").Value = 40.706319702602
$ws.Range("M24").Value = 5.285118219749

# Row 25
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = -21.428571428571
$ws.Range("G25").Value = 38
$ws.Range("H25").Value = 13.157894736842
$ws.Range("I25").Value = 448
$ws.Range("J25").Value = 367
$ws.Range("K25").Value = 22.070844686648
$ws.Range("L25").Value = 44.516129032258
$ws.Range("M25").Value = -39.130434782608

# Row 26
$ws.Range("G26").Value = 3
$ws.Range("J26").Value = 27
$ws.Range("K26").Value = -7.407407407407

# Row 27
$ws.Range("D27").Value = 3
$ws.Range("I14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -75
$ws.Range("J27").Value = 31
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = -35.416666666666

# Row 28
$ws.Range("D28").Value = 1
$ws.Range("I14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("F28").NumberFormat = "@"
$ws.Range("F28").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F28").PasteSpecial(-4122)
$ws.Range("G28").Value = 1
$ws.Range("I14").Copy()
$ws.Range("G28").PasteSpecial(-4122)
$ws.Range("H28").Value = -100
$ws.Range("K14").Copy()
$ws.Range("H28").PasteSpecial(-4122)
$ws.Range("J28").Value = 35
$ws.Range("K28").Value = -57.142857142857
$ws.Range("L28").Value = -61.538461538461
$ws.Range("M28").Value = -72.222222222222

# Row 29
$ws.Range("D29").Value = 1
$ws.Range("I14").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E29").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("F29").NumberFormat = "@"
$ws.Range("F29").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("G29").Value = 1
$ws.Range("I14").Copy()
$ws.Range("G29").PasteSpecial(-4122)
$ws.Range("H29").Value = -100
$ws.Range("K14").Copy()
$ws.Range("H29").PasteSpecial(-4122)
$ws.Range("J29").Value = 27
$ws.Range("K29").Value = -55.555555555555
$ws.Range("L29").Value = -60
$ws.Range("M29").Value = -73.91304347826

$excel.CutCopyMode = $false